# Update cryptocurrency price/volume snapshot (GitHub Actions refresh).
# Numeric-looking "Price" text values are written with a leading single-quote
# (Excel's text-quote prefix) so they stay literal text instead of being
# auto-coerced into numbers - matching the sheet's original inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '30.141.11'
$ws.Range('E2').Value = '  +0.10%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '1.915.38'
$ws.Range('E3').Value = '  +0.28%  '

# Row 4: TetherUSD
$ws.Range('D4').Value = '''1.000'
$ws.Range('E4').Value = '  -0.05%  '

# Row 5: XRP
$ws.Range('D5').Value = '''0.7945'
$ws.Range('E5').Value = '  +7.25%  '

# Row 6: BNB
$ws.Range('D6').Value = '''243.00'
$ws.Range('E6').Value = '  -0.61%  '

# Row 7: USDC
$ws.Range('D7').Value = '''1.000'
$ws.Range('E7').Value = '  -0.03%  '

# Row 8: Cardano
$ws.Range('D8').Value = '''0.3189'
$ws.Range('E8').Value = '  +3.05%  '

# Row 9: Solana
$ws.Range('D9').Value = '''26.40'
$ws.Range('E9').Value = '  -0.39%  '

# Row 10: Dogecoin
$ws.Range('D10').Value = '''0.06971'
$ws.Range('E10').Value = '  -0.14%  '

# Row 11: TRON
$ws.Range('D11').Value = '''0.08021'
$ws.Range('E11').Value = '  -0.72%  '

# Row 12: Polygon
$ws.Range('E12').Value = '  -2.28%  '

# Row 13: WrappedEther
$ws.Range('D13').Value = '1.914.24'
$ws.Range('E13').Value = '  +0.33%  '

# Row 14: Polkadot
$ws.Range('D14').Value = '''5.233'
$ws.Range('E14').Value = '  -1.89%  '

# Row 15: Litecoin
$ws.Range('D15').Value = '''93.66'
$ws.Range('E15').Value = '  +1.44%  '

# Row 16: WrappedBTC
$ws.Range('D16').Value = '30.156.33'
$ws.Range('E16').Value = '  +0.14%  '

# Row 17: Avalanche
$ws.Range('D17').Value = '''14.07'
$ws.Range('E17').Value = '  -1.75%  '

# Row 18: Uniswap
$ws.Range('D18').Value = '''6.011'
$ws.Range('E18').Value = '  -1.00%  '

# Row 19: BitcoinCash
$ws.Range('D19').Value = '''249.51'
$ws.Range('E19').Value = '  +3.77%  '

# Row 20: ShibaInu
$ws.Range('D20').Value = '''0.000007837'
$ws.Range('E20').Value = '  -0.13%  '

# Row 21: Dai
$ws.Range('D21').Value = '''0.9998'
$ws.Range('E21').Value = '  -0.10%  '

# Row 22: WrappedliquidstakedEther2.0
$ws.Range('D22').Value = '2.149.00'
$ws.Range('E22').Value = '  -0.75%  '

# Row 23: BinanceUSD
$ws.Range('D23').Value = '''1.000'
$ws.Range('E23').Value = '  -0.03%  '

# Row 24: Chainlink
$ws.Range('E24').Value = '  -1.11%  '

# Row 25: Monero
$ws.Range('D25').Value = '''169.23'
$ws.Range('E25').Value = '  +1.21%  '

# Row 26: Cosmos
$ws.Range('D26').Value = '''9.348'
$ws.Range('E26').Value = '  -0.61%  '

# Row 27: Stellar
$ws.Range('D27').Value = '''0.1398'
$ws.Range('E27').Value = '  +9.32%  '

# Row 28: EthereumClassic
$ws.Range('D28').Value = '''19.02'
$ws.Range('E28').Value = '  +0.19%  '

# Row 29: LidoDAOToken
$ws.Range('D29').Value = '''2.058'
$ws.Range('E29').Value = '  -0.09%  '

# Row 30: Toncoin
$ws.Range('E30').Value = '  +2.59%  '

# Row 31: PancakeSwap
$ws.Range('D31').Value = '''1.528'
$ws.Range('E31').Value = '  -1.43%  '

# Row 32: Filecoin
$ws.Range('D32').Value = '''4.374'
$ws.Range('E32').Value = '  +0.73%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range('D33').Value = '''4.126'
$ws.Range('E33').Value = '  +0.99%  '

# Row 34: Hedera
$ws.Range('D34').Value = '''0.05389'
$ws.Range('E34').Value = '  +4.52%  '

# Row 35: ARBITRUM
$ws.Range('D35').Value = '''1.272'
$ws.Range('E35').Value = '  -2.85%  '

# Row 36: ImmutableX
$ws.Range('D36').Value = '''0.7403'
$ws.Range('E36').Value = '  -1.13%  '

# Row 37: HuobiToken
$ws.Range('D37').Value = '''2.730'
$ws.Range('E37').Value = '  +0.18%  '

# Row 38: VeChain
$ws.Range('D38').Value = '''0.01933'
$ws.Range('E38').Value = '  -1.27%  '

# Row 39: MXToken
$ws.Range('D39').Value = '''2.796'
$ws.Range('E39').Value = '  -0.02%  '

# Row 40: TheSandbox
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '''6.205'
$ws.Range('E40').Value = '  -2.18%  '

# Row 41: FraxShare
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '''0.4471'
$ws.Range('E41').Value = '  -0.85%  '

# Row 42: Aave
$ws.Range('D42').Value = '''72.77'
$ws.Range('E42').Value = '  -2.34%  '

# Row 43: RenderToken
$ws.Range('D43').Value = '''1.907'
$ws.Range('E43').Value = '  -3.78%  '

# Row 44: PaxDollar
$ws.Range('D44').Value = '''1.000'
$ws.Range('E44').Value = '  -0.15%  '

# Row 45: TrustWalletToken
$ws.Range('D45').Value = '''0.8354'
$ws.Range('E45').Value = '  -0.64%  '

# Row 46: Aptos
$ws.Range('D46').Value = '''7.628'
$ws.Range('E46').Value = '  -1.60%  '

# Row 47: Quant
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '''9.896'
$ws.Range('E47').Value = '  -0.56%  '

# Row 48: EnergySwap
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '''100.71'
$ws.Range('E48').Value = '  -1.20%  '

# Row 49: RocketPoolETH
$ws.Range('D49').Value = '2.059.12'
$ws.Range('E49').Value = '  -0.56%  '

# Row 50: Elrond
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '''967.34'
$ws.Range('E50').Value = '  +4.54%  '

# Row 51: Maker
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').Value = '''36.55'
$ws.Range('E51').Value = '  -0.61%  '
